# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 13630
$wsExpo.Range("F3").Value = 321
$wsExpo.Range("F4").Value = 660
$wsExpo.Range("F5").Value = 227
$wsExpo.Range("F6").Value = 476
$wsExpo.Range("F7").Value = 1383

# Sheet "全部类型" (all types combined)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 13630
$wsAll.Range("F3").Value = 321
$wsAll.Range("F4").Value = 660
$wsAll.Range("F5").Value = 227
$wsAll.Range("F8").Value = 476
$wsAll.Range("F9").Value = 1383
